$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove buybacks_cash_returned column (T) data for rows 2-3
$ws.Range("T2:T3").ClearContents()

# Row 2 updates
$ws.Range("D2").Value = -0.0532
$ws.Range("G2").Value = 0.02504862953138816
$ws.Range("H2").Value = 0.002060123784261716
$ws.Range("I2").Value = -0.01282051282051282
$ws.Range("J2").Value = -0.01282051282051282
$ws.Range("K2").Value = -2.62
$ws.Range("L2").Value = -0.02316534040671972
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = -0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = -0
$ws.Range("U2").Value = 13.7
$ws.Range("V2").Value = 0.8726114649681529
$ws.Range("W2").Value = -0.1578313253012048
$ws.Range("X2").Value = 0.1072654182157499
$ws.Range("Y2").Value = -0.2650967435169547
$ws.Range("Z2").Value = 4.491660047656869
$ws.Range("AA2").Value = -0.05758538522637012
$ws.Range("AB2").Value = 0.05992321584424368
$ws.Range("AC2").Value = -0.1175086010706138
$ws.Range("AD2").Value = 17.8
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 17.8
$ws.Range("AG2").Value = 4.100000000000001
$ws.Range("AH2").Value = 0.5313432835820896
$ws.Range("AI2").Value = 0.5345345345345346
$ws.Range("AJ2").Value = 0.2070707070707071
$ws.Range("AK2").Value = 0.2091836734693878
$ws.Range("AL2").Value = 1.05
$ws.Range("AM2").Value = 1.044
$ws.Range("AN2").Value = 10.17142857142857
$ws.Range("AO2").Value = -1.380952380952381
$ws.Range("AP2").Value = 2.342857142857144
$ws.Range("AQ2").Value = -1.388888888888889

# Row 3 updates
$ws.Range("D3").Value = -0.0532
$ws.Range("G3").Value = 0.02504862953138816
$ws.Range("H3").Value = 0.002060123784261716
$ws.Range("I3").Value = -0.01282051282051282
$ws.Range("J3").Value = -0.01282051282051282
$ws.Range("K3").Value = -2.62
$ws.Range("L3").Value = -0.02316534040671972
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = 0
$ws.Range("U3").Value = 13.7
$ws.Range("V3").Value = 0.8726114649681529
$ws.Range("W3").Value = -0.1578313253012048
$ws.Range("X3").Value = 0.1072654182157499
$ws.Range("Y3").Value = -0.2650967435169547
$ws.Range("Z3").Value = 4.491660047656869
$ws.Range("AA3").Value = -0.05758538522637012
$ws.Range("AB3").Value = 0.05992321584424368
$ws.Range("AC3").Value = -0.1175086010706138
$ws.Range("AD3").Value = 17.8
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 17.8
$ws.Range("AG3").Value = 4.100000000000001
$ws.Range("AH3").Value = 0.5313432835820896
$ws.Range("AI3").Value = 0.5345345345345346
$ws.Range("AJ3").Value = 0.2070707070707071
$ws.Range("AK3").Value = 0.2091836734693878
$ws.Range("AL3").Value = 1.05
$ws.Range("AM3").Value = 1.044
$ws.Range("AN3").Value = 10.17142857142857
$ws.Range("AO3").Value = -1.380952380952381
$ws.Range("AP3").Value = 2.342857142857144
$ws.Range("AQ3").Value = -1.388888888888889
